{"js": "// Update the date line and the 25 \"three-digit \u00f7 one-digit\" answer cells\n// to the new values described in the commit diff. Each old value is\n// unique within the document, so a scoped search-and-replace (in\n// \"Replace\" mode, which preserves the existing run formatting) is used\n// for every pair.\nconst replacements = [\n  [\"2025-04-27 Sunday\", \"2025-04-28 Monday\"],\n  [\"109\u00f77=15, 4\", \"346\u00f78=43, 2\"],\n  [\"323\u00f72=161, 1\", \"293\u00f79=32, 5\"],\n  [\"186\u00f78=23, 2\", \"347\u00f79=38, 5\"],\n  [\"158\u00f77=22, 4\", \"417\u00f74=104, 1\"],\n  [\"902\u00f72=451, 0\", \"266\u00f79=29, 5\"],\n  [\"577\u00f75=115, 2\", \"600\u00f79=66, 6\"],\n  [\"827\u00f76=137, 5\", \"157\u00f79=17, 4\"],\n  [\"931\u00f79=103, 4\", \"177\u00f72=88, 1\"],\n  [\"501\u00f78=62, 5\", \"528\u00f73=176, 0\"],\n  [\"230\u00f73=76, 2\", \"986\u00f73=328, 2\"],\n  [\"363\u00f72=181, 1\", \"884\u00f72=442, 0\"],\n  [\"427\u00f72=213, 1\", \"870\u00f75=174, 0\"],\n  [\"101\u00f76=16, 5\", \"401\u00f77=57, 2\"],\n  [\"928\u00f76=154, 4\", \"598\u00f77=85, 3\"],\n  [\"296\u00f79=32, 8\", \"758\u00f77=108, 2\"],\n  [\"669\u00f74=167, 1\", \"841\u00f75=168, 1\"],\n  [\"607\u00f78=75, 7\", \"453\u00f74=113, 1\"],\n  [\"320\u00f73=106, 2\", \"612\u00f77=87, 3\"],\n  [\"701\u00f78=87, 5\", \"446\u00f74=111, 2\"],\n  [\"674\u00f76=112, 2\", \"128\u00f72=64, 0\"],\n  [\"931\u00f77=133, 0\", \"845\u00f74=211, 1\"],\n  [\"442\u00f73=147, 1\", \"666\u00f75=133, 1\"],\n  [\"325\u00f72=162, 1\", \"701\u00f77=100, 1\"],\n  [\"998\u00f78=124, 6\", \"781\u00f73=260, 1\"],\n  [\"994\u00f73=331, 1\", \"696\u00f74=174, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 \"three-digit \u00f7 one-digit\" answer cells\n# to the new values described in the commit diff. Each old value is\n# unique within the document, so a single Find/Replace (ReplaceOne) pass\n# per pair is sufficient and keeps the existing run formatting intact.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-04-27 Sunday\", \"2025-04-28 Monday\"),\n    @(\"109\u00f77=15, 4\", \"346\u00f78=43, 2\"),\n    @(\"323\u00f72=161, 1\", \"293\u00f79=32, 5\"),\n    @(\"186\u00f78=23, 2\", \"347\u00f79=38, 5\"),\n    @(\"158\u00f77=22, 4\", \"417\u00f74=104, 1\"),\n    @(\"902\u00f72=451, 0\", \"266\u00f79=29, 5\"),\n    @(\"577\u00f75=115, 2\", \"600\u00f79=66, 6\"),\n    @(\"827\u00f76=137, 5\", \"157\u00f79=17, 4\"),\n    @(\"931\u00f79=103, 4\", \"177\u00f72=88, 1\"),\n    @(\"501\u00f78=62, 5\", \"528\u00f73=176, 0\"),\n    @(\"230\u00f73=76, 2\", \"986\u00f73=328, 2\"),\n    @(\"363\u00f72=181, 1\", \"884\u00f72=442, 0\"),\n    @(\"427\u00f72=213, 1\", \"870\u00f75=174, 0\"),\n    @(\"101\u00f76=16, 5\", \"401\u00f77=57, 2\"),\n    @(\"928\u00f76=154, 4\", \"598\u00f77=85, 3\"),\n    @(\"296\u00f79=32, 8\", \"758\u00f77=108, 2\"),\n    @(\"669\u00f74=167, 1\", \"841\u00f75=168, 1\"),\n    @(\"607\u00f78=75, 7\", \"453\u00f74=113, 1\"),\n    @(\"320\u00f73=106, 2\", \"612\u00f77=87, 3\"),\n    @(\"701\u00f78=87, 5\", \"446\u00f74=111, 2\"),\n    @(\"674\u00f76=112, 2\", \"128\u00f72=64, 0\"),\n    @(\"931\u00f77=133, 0\", \"845\u00f74=211, 1\"),\n    @(\"442\u00f73=147, 1\", \"666\u00f75=133, 1\"),\n    @(\"325\u00f72=162, 1\", \"701\u00f77=100, 1\"),\n    @(\"998\u00f78=124, 6\", \"781\u00f73=260, 1\"),\n    @(\"994\u00f73=331, 1\", \"696\u00f74=174, 0\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
